# Se corrigen los formatos de fecha y numero
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("B2").Value = 12300000
$ws.Range("F2").Value = 45523
$ws.Range("F2").NumberFormat = "yyyy-mm-dd"
$ws.Range("F2").NumberFormat = "YYYY-MM-DD"
$ws.Range("G2").Value = 45524
$ws.Range("G2").NumberFormat = "YYYY-MM-DD"
$ws.Range("H2").Value = 11247.5571
$ws.Range("I2").Value = 1093.5708

# --- Row 3 ---
$ws.Range("B3").Value = 5323478
$ws.Range("F3").Value = 45523
$ws.Range("F3").NumberFormat = "YYYY-MM-DD"
$ws.Range("G3").Value = 45524
$ws.Range("G3").NumberFormat = "YYYY-MM-DD"
$ws.Range("H3").Value = 4867.9773
$ws.Range("I3").Value = 1093.5708

# --- Row 4 ---
$ws.Range("B4").Value = 34151465
$ws.Range("F4").Value = 45523
$ws.Range("F4").NumberFormat = "YYYY-MM-DD"
$ws.Range("G4").Value = 45524
$ws.Range("G4").NumberFormat = "YYYY-MM-DD"
$ws.Range("H4").Value = 31229.3133
$ws.Range("I4").Value = 1093.5708

# --- Row 5 ---
$ws.Range("B5").Value = 15002401
$ws.Range("F5").Value = 45523
$ws.Range("F5").NumberFormat = "YYYY-MM-DD"
$ws.Range("G5").Value = 45524
$ws.Range("G5").NumberFormat = "YYYY-MM-DD"
$ws.Range("H5").Value = 13718.7284
$ws.Range("I5").Value = 1093.5708
